$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.826.35"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.636.70"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D5").Value = "'215.85"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'0.5071"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.2579"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.06436"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Value = "'19.61"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "'4.280"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "1.863.46"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "1.635.12"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "'0.5642"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "0.0₅7608"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "'63.21"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "25.851.63"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'195.65"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  -2.94%  "
$ws.Range("D22").Value = "'9.891"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "'6.097"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'1.792"
$ws.Range("E25").Value = "  -5.46%  "
$ws.Range("D26").Value = "'0.1271"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Value = "'139.76"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").Value = "'6.791"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").Value = "'15.51"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").Value = "'1.242"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'0.04890"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "'3.301"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").Value = "'3.227"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").Value = "'1.559"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'0.9041"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "'2.577"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "1.129.93"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "'0.5511"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "'0.01563"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'0.9941"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "'5.531"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "'0.8003"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "'97.75"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").Value = "1.773.54"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").Value = "'55.41"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05051"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.679"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("E51").Value = "  +0.07%  "
